$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 44: only columns A and C are populated (same shape as the other
# blank "placeholder" day-rows), C uses the existing short-date style (s="4")
# that is already used by C43, so copy formatting from there first.
$ws.Cells.Item(43, 3).Copy() | Out-Null
$ws.Cells.Item(44, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 3).Value = 45319

# --- Row 45: same shape as row 44.
$ws.Cells.Item(43, 3).Copy() | Out-Null
$ws.Cells.Item(45, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 3).Value = 45320

# --- Row 46: a full data row. Copy the style/number-format layout from row
# 43 (the prior full data row) across A:N, then fill in the values/formulas.
$ws.Range("A43:N43").Copy() | Out-Null
$ws.Range("A46:N46").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 35
$ws.Cells.Item(46, 3).Value = 45321
$ws.Cells.Item(46, 4).Value = 1.157
$ws.Cells.Item(46, 5).Value = 1
$ws.Cells.Item(46, 6).Value = 1015
$ws.Cells.Item(46, 7).Formula = '=F46*E46*D46'
$ws.Cells.Item(46, 8).Formula = '=I43'
$ws.Cells.Item(46, 9).Formula = '=H46+G46-F46'
$ws.Cells.Item(46, 10).Formula = '=I46-H46'
$ws.Cells.Item(46, 11).Formula = '=I46/$H$2-1'
$ws.Cells.Item(46, 12).Value = "ASIA"
$ws.Cells.Item(46, 13).Value = "ESPORTS"
$ws.Cells.Item(46, 14).Value = "LCK"
# Column G has no explicit style on the target row (General format); the
# formula above can pick up D46's style via its cell reference, so force it
# back to the default/Normal style to match.
$ws.Cells.Item(46, 7).Style = "Normal"

# --- Row 47: another full data row, same layout as row 46.
$ws.Range("A43:N43").Copy() | Out-Null
$ws.Range("A47:N47").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = 36
$ws.Cells.Item(47, 3).Value = 45321
$ws.Cells.Item(47, 4).Value = 1.115
$ws.Cells.Item(47, 5).Value = 1
$ws.Cells.Item(47, 6).Value = 500
$ws.Cells.Item(47, 7).Formula = '=F47*E47*D47'
$ws.Cells.Item(47, 8).Formula = '=I46'
$ws.Cells.Item(47, 9).Formula = '=H47+G47-F47'
$ws.Cells.Item(47, 10).Formula = '=I47-H47'
$ws.Cells.Item(47, 11).Formula = '=I47/$H$2-1'
$ws.Cells.Item(47, 12).Value = "ASIA"
$ws.Cells.Item(47, 13).Value = "TABLE TENNIS"
$ws.Cells.Item(47, 14).Value = "LIGA PRO"
$ws.Cells.Item(47, 7).Style = "Normal"

$excel.CutCopyMode = 0

# --- View state: the author scrolled/selected a different area after
# adding the new rows.
$ws.Application.GoTo($ws.Range("A21"), $true)
$ws.Range("O54").Select() | Out-Null
